$wb = $excel.ActiveWorkbook

# --- Overview sheet: update status for the a8180a71 row (row 3) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Handback transform failed"
$wsOverview.Range("F3").Value = "Handback transform failed"

# --- zh-cn sheet: update Status + Error Detail for row 3, widen Error Detail column ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Handback transform failed"
$wsZhCn.Range("P3").Value = "Handback file name: ujdmj2me.3jy is different with handoff file name: a8180a71-a44a-4e23-89f6-6e2ceefde047.0be0af7b565512c22467dae631a41eb038fd06f6.zh-cn."
$wsZhCn.Columns.Item(16).ColumnWidth = 39.166666666666664

# --- de-de sheet: update Status + Error Detail for row 3, widen Error Detail column ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Handback transform failed"
$wsDeDe.Range("P3").Value = "Handback file name: ujdmj2me.3jy is different with handoff file name: a8180a71-a44a-4e23-89f6-6e2ceefde047.0be0af7b565512c22467dae631a41eb038fd06f6.de-de."
$wsDeDe.Columns.Item(16).ColumnWidth = 39.166666666666664
